$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1) Remove the CLARIBEL PINILLA AMAYA (CC / 1007314918) block of 4 rows entirely.
#    This shifts the DIDIER LOUIS JOSEPH rows (and everything below) up by 4.
$ws.Rows("16:19").Delete()

# 2) Update the "VALOR MORA" total (E11): 1841334 -> 1553334
$ws.Range("E11").Value = 1553334

# 3) Update "Cant. Trabajadores" (C13): 2 -> 1 (only one worker left in the report)
$ws.Range("C13").Value = 1

# 4) Re-order the DIDIER LOUIS JOSEPH period rows (now rows 16-23) from
#    descending (2310..2209) to ascending (2209..2310), and move the
#    153334 "Valor Mora" amount from the first row to the last row,
#    matching the new "parte 1" statement ordering.
$periods = @("2209", "2304", "2305", "2306", "2307", "2308", "2309", "2310")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    if ($row -eq 23) {
        $ws.Range("F$row").Value = 153334
    } else {
        $ws.Range("F$row").Value = 200000
    }
    $ws.Range("G$row").Value = 5000000
}
